# Weekly update: a new price record for "Feria Lagunitas de Puerto Montt - Espinaca"
# is inserted as the new row 33, pushing all the previously existing rows
# (old rows 33..62) down by one (new rows 34..63).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 33; this shifts rows 33:62 down to 34:63
# and extends the sheet dimension accordingly.
$ws.Rows("33:33").Insert()

# Populate the newly inserted row 33 with the new weekly record.
$ws.Range("A33").Value = 4
$ws.Range("B33").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C33").Value = "Los Lagos"
$ws.Range("D33").Value = 45072
$ws.Range("E33").Value = 10
$ws.Range("F33").Value = 100112012
$ws.Range("G33").Value = "Espinaca"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 25
$ws.Range("K33").Value = 13000
$ws.Range("L33").Value = 13000
$ws.Range("M33").Value = 13000
$ws.Range("N33").Value = "$/cuna 10 kilos"
$ws.Range("O33").Value = "Región Metropolitana"
$ws.Range("P33").Value = 1300
$ws.Range("Q33").Value = 10
$ws.Range("R33").Value = "Hortaliza"
